$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.673.97"
$ws.Range("E2").Value = "  +0.14%  "
$ws.Range("D3").Value = "1.597.36"
$ws.Range("E3").Value = "  -0.02%  "
$ws.Range("E4").Value = "  +0.17%  "
$ws.Range("D5").Value = "'211.47"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.34%  "
$ws.Range("E6").Value = "  +0.01%  "
$ws.Range("E7").Value = "  +0.21%  "
$ws.Range("E8").Value = "  +0.20%  "
$ws.Range("E9").Value = "  +0.71%  "
$ws.Range("D10").Value = "'19.46"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.47%  "
$ws.Range("E11").Value = "  +0.25%  "
$ws.Range("D12").Value = "1.821.16"
$ws.Range("E12").Value = "  -0.07%  "
$ws.Range("D13").Value = "1.616.30"
$ws.Range("E13").Value = "  +1.43%  "
$ws.Range("D14").Value = "'4.02"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.18%  "
$ws.Range("D15").Value = "'0.523"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.56%  "
$ws.Range("D16").Value = "'65.04"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.58%  "
$ws.Range("D17").Value = "26.643.04"
$ws.Range("E17").Value = "  +0.08%  "
$ws.Range("D18").Value = "0.0₃0738"
$ws.Range("E18").Value = "  +1.25%  "
$ws.Range("E19").Value = "  +0.17%  "
$ws.Range("E20").Value = "  +0.14%  "
$ws.Range("D21").Value = "'7.05"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +4.62%  "
$ws.Range("D22").Value = "'4.29"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.87%  "
$ws.Range("E23").Value = "  +2.15%  "
$ws.Range("E24").Value = "  +1.12%  "
$ws.Range("D25").Value = "'143.46"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.46%  "
$ws.Range("E26").Value = "  +0.16%  "
$ws.Range("D27").Value = "'7.12"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.50%  "
$ws.Range("E28").Value = "  -0.87%  "
$ws.Range("D29").Value = "'15.30"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.23%  "
$ws.Range("E30").Value = "  +1.69%  "
$ws.Range("E31").Value = "  +0.06%  "
$ws.Range("E32").Value = "  +0.12%  "
$ws.Range("E33").Value = "  +0.50%  "
$ws.Range("D34").Value = "1.288.55"
$ws.Range("E34").Value = "  -0.19%  "
$ws.Range("D35").Value = "'0.619"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -5.30%  "
$ws.Range("E36").Value = "  +0.57%  "
$ws.Range("E37").Value = "  -0.15%  "
$ws.Range("D38").Value = "'0.0170"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.35%  "
$ws.Range("E39").Value = "  -1.65%  "
$ws.Range("E40").Value = "  +15.58%  "
$ws.Range("E41").Value = "  +1.06%  "
$ws.Range("D42").Value = "'2.19"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.55%  "
$ws.Range("D43").Value = "'0.781"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Value = "'63.21"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.96%  "
$ws.Range("D45").Value = "1.732.53"
$ws.Range("D46").Value = "'91.06"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.14%  "
$ws.Range("E47").Value = "  -2.89%  "
$ws.Range("E48").Value = "  +1.19%  "
$ws.Range("E49").Value = "  +0.90%  "
$ws.Range("E50").Value = "  +0.46%  "
$ws.Range("D51").Value = "'7.37"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.94%  "
